$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- JSON payloads for the value_json column (rows 2..5) ---

$json2 = @"
[
	{
		"code": "101",
		"value": "A",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "A+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "A-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "104",
		"value": "B",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "105",
		"value": "B+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "106",
		"value": "B-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "107",
		"value": "AB",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "108",
		"value": "AB+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "109",
		"value": "AB-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "110",
		"value": "O",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "111",
		"value": "O+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "112",
		"value": "O-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "113",
		"value": "Don't Know",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "114",
		"value": "Not Applicable",
		"langCode": "eng",
		"active": true
	}
]
"@

$json3 = @"
[
	{
		"code": "101",
		"value": "Single",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "Married",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "Widowed",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "104",
		"value": "Divorced",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "105",
		"value": "Legally Separated",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "106",
		"value": "Annulled",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "107",
		"value": "Nullified",
		"langCode": "eng",
		"active": true
	}
]
"@

$json4 = @"
[
	{
		"code": "Document-based",
		"value": "Document-based",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "Introducer-based",
		"value": "Introducer-based",
		"langCode": "eng",
		"active": true
	}
]
"@

$json5 = @"
[
	{
		"code": "101",
		"value": "Pick-up",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "Delivery to permanent address",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "Delivery to present address",
		"langCode": "eng",
		"active": true
	}
]
"@

# --- Propagate the bold/bordered header style (currently on A1:G1) to the
#     new header cell H1, and to the new id-index column cells A2:A5,
#     before any values move/overwrite. ---
$ws.Range("A1").Copy($ws.Range("H1"))
$ws.Range("A1").Copy($ws.Range("A2:A5"))

# --- Row 1 headers: insert lang_code before id, rename dataType -> data_type,
#     shift everything right and drop the old A1 (id moves to C1). ---
$ws.Range("A1").Clear()
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "id"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "description"
$ws.Range("F1").Value = "data_type"
$ws.Range("G1").Value = "value_json"
$ws.Range("H1").Value = "is_active"

# --- Row 2: existing bloodType record, now with a leading numeric id (0)
#     and lang_code column, description/value_json rewritten. ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = 10001
$ws.Range("D2").Value = "bloodType"
$ws.Range("E2").Value = "Blood Type"
$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = $json2
$ws.Range("H2").Value = $true
$ws.Rows(2).AutoFit()

# --- Row 3: maritalStatus (new) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eng"
$ws.Range("C3").Value = 10002
$ws.Range("D3").Value = "maritalStatus"
$ws.Range("E3").Value = "Marital Status"
$ws.Range("F3").Value = "string"
$ws.Range("G3").Value = $json3
$ws.Range("H3").Value = $true
$ws.Rows(3).AutoFit()

# --- Row 4: registrationType (new) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "eng"
$ws.Range("C4").Value = 10003
$ws.Range("D4").Value = "registrationType"
$ws.Range("E4").Value = "Registration Type"
$ws.Range("F4").Value = "string"
$ws.Range("G4").Value = $json4
$ws.Range("H4").Value = $true
$ws.Rows(4).AutoFit()

# --- Row 5: modeOfClaim (new) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "eng"
$ws.Range("C5").Value = 10004
$ws.Range("D5").Value = "modeOfClaim"
$ws.Range("E5").Value = "Mode of Claim"
$ws.Range("F5").Value = "string"
$ws.Range("G5").Value = $json5
$ws.Range("H5").Value = $true
$ws.Rows(5).AutoFit()
